$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = -2.288097732945488
$ws.Cells.Item(2, 4).Value = -4.339917659759537
$ws.Cells.Item(2, 5).Value = -3.936509232772027
$ws.Cells.Item(2, 6).Value = -0.2191115766763687
$ws.Cells.Item(2, 7).Value = -0.8925584554672241
$ws.Cells.Item(2, 8).Value = 0.058589544147253

$ws.Cells.Item(3, 3).Value = -3.323771476745605
$ws.Cells.Item(3, 4).Value = -0.8808293342590332
$ws.Cells.Item(3, 5).Value = -3.401906013488769
$ws.Cells.Item(3, 6).Value = -0.4737097918987274
$ws.Cells.Item(3, 7).Value = -0.3764378130435943
$ws.Cells.Item(3, 8).Value = -0.0346210934221744

$ws.Cells.Item(4, 3).Value = -5.066901056390067
$ws.Cells.Item(4, 4).Value = 1.504574901179274
$ws.Cells.Item(4, 5).Value = -3.703831045251145
$ws.Cells.Item(4, 6).Value = 2.269346237182617
$ws.Cells.Item(4, 7).Value = 2.84225869178772
$ws.Cells.Item(4, 8).Value = -1.580319762229919

$ws.Cells.Item(5, 3).Value = -4.698070149672659
$ws.Cells.Item(5, 4).Value = -0.4452685556913691
$ws.Cells.Item(5, 5).Value = -8.655812690132544
$ws.Cells.Item(5, 6).Value = 2.80038046836853
$ws.Cells.Item(5, 7).Value = -2.297109603881836
$ws.Cells.Item(5, 8).Value = -0.6647250056266785

$ws.Cells.Item(6, 3).Value = -6.691070581737309
$ws.Cells.Item(6, 4).Value = -6.563599611583482
$ws.Cells.Item(6, 5).Value = -8.057312513652606
$ws.Cells.Item(6, 6).Value = -2.487059593200684
$ws.Cells.Item(6, 7).Value = -0.9021458625793456
$ws.Cells.Item(6, 8).Value = 2.389388084411621

$ws.Cells.Item(7, 3).Value = -4.573285604778092
$ws.Cells.Item(7, 4).Value = -17.30385228207236
$ws.Cells.Item(7, 5).Value = 12.64209365844724
$ws.Cells.Item(7, 6).Value = -9.80595874786377
$ws.Cells.Item(7, 7).Value = 2.084323167800904
$ws.Cells.Item(7, 8).Value = 0.3504720032215118

$ws.Cells.Item(8, 3).Value = -1.509485269847658
$ws.Cells.Item(8, 4).Value = -5.039251478094737
$ws.Cells.Item(8, 5).Value = 6.722280000385433
$ws.Cells.Item(8, 6).Value = 1.612078070640564
$ws.Cells.Item(8, 7).Value = 6.0146164894104
$ws.Cells.Item(8, 8).Value = 2.238120555877685

$ws.Cells.Item(9, 3).Value = -1.532178590172217
$ws.Cells.Item(9, 4).Value = -16.01636389682171
$ws.Cells.Item(9, 5).Value = 6.974593714663866
$ws.Cells.Item(9, 6).Value = 3.856590270996094
$ws.Cells.Item(9, 7).Value = -5.252686023712158
$ws.Cells.Item(9, 8).Value = 0.6668555736541748

$ws.Cells.Item(10, 3).Value = -22.63299742497886
$ws.Cells.Item(10, 4).Value = -9.195889523154877
$ws.Cells.Item(10, 5).Value = -11.14407699986499
$ws.Cells.Item(10, 6).Value = 7.994077205657959
$ws.Cells.Item(10, 7).Value = 4.401273250579834
$ws.Cells.Item(10, 8).Value = 0.2966761589050293

$ws.Cells.Item(11, 3).Value = -14.47936531117094
$ws.Cells.Item(11, 4).Value = -0.6024172117835604
$ws.Cells.Item(11, 5).Value = -11.01582479476932
$ws.Cells.Item(11, 6).Value = 0.5494767427444458
$ws.Cells.Item(11, 7).Value = 2.293114900588989
$ws.Cells.Item(11, 8).Value = -1.938781261444092

$ws.Cells.Item(12, 3).Value = -1.626747369766244
$ws.Cells.Item(12, 4).Value = -8.822486783328847
$ws.Cells.Item(12, 5).Value = -5.710294397253755
$ws.Cells.Item(12, 6).Value = -5.656754016876221
$ws.Cells.Item(12, 7).Value = -0.6400907635688782
$ws.Cells.Item(12, 8).Value = 0.8974853157997131

$ws.Cells.Item(13, 3).Value = 0.4268563170182262
$ws.Cells.Item(13, 4).Value = -7.612172829477366
$ws.Cells.Item(13, 5).Value = 6.265647787796897
$ws.Cells.Item(13, 6).Value = -7.037869453430176
$ws.Cells.Item(13, 7).Value = -3.455584764480591
$ws.Cells.Item(13, 8).Value = -5.363073825836182

$ws.Cells.Item(14, 3).Value = 0.09527274181970302
$ws.Cells.Item(14, 4).Value = 5.863329686616582
$ws.Cells.Item(14, 5).Value = 12.21953241448658
$ws.Cells.Item(14, 6).Value = -0.4055328667163849
$ws.Cells.Item(14, 7).Value = 9.243432998657228
$ws.Cells.Item(14, 8).Value = -5.18677282333374

$ws.Cells.Item(15, 3).Value = 8.31707101119181
$ws.Cells.Item(15, 4).Value = -16.11775568911882
$ws.Cells.Item(15, 5).Value = 2.399655392295399
$ws.Cells.Item(15, 6).Value = 7.240935325622559
$ws.Cells.Item(15, 7).Value = -3.536012172698975
$ws.Cells.Item(15, 8).Value = -0.926780104637146

$ws.Cells.Item(16, 3).Value = 23.51702158074626
$ws.Cells.Item(16, 4).Value = -39.05408849214243
$ws.Cells.Item(16, 5).Value = -0.4082389630769931
$ws.Cells.Item(16, 6).Value = 1.342033505439758
$ws.Cells.Item(16, 7).Value = 2.555702686309814
$ws.Cells.Item(16, 8).Value = -1.59523355960846

$ws.Cells.Item(17, 3).Value = -2.529885241859907
$ws.Cells.Item(17, 4).Value = 6.300281524658134
$ws.Cells.Item(17, 5).Value = -23.18308167708546
$ws.Cells.Item(17, 6).Value = -3.227951049804688
$ws.Cells.Item(17, 7).Value = 2.566887855529785
$ws.Cells.Item(17, 8).Value = 1.657018899917602

$ws.Cells.Item(18, 3).Value = -4.514081553409033
$ws.Cells.Item(18, 4).Value = -21.8527907321327
$ws.Cells.Item(18, 5).Value = -10.32276263989904
$ws.Cells.Item(18, 6).Value = -7.749465942382812
$ws.Cells.Item(18, 7).Value = -4.05479621887207
$ws.Cells.Item(18, 8).Value = 8.297345161437988

$ws.Cells.Item(19, 3).Value = 38.92504751054884
$ws.Cells.Item(19, 4).Value = -55.46849928404185
$ws.Cells.Item(19, 5).Value = 19.13083159296148
$ws.Cells.Item(19, 6).Value = 5.375124931335449
$ws.Cells.Item(19, 7).Value = -9.55262565612793
$ws.Cells.Item(19, 8).Value = -2.865561485290528

$ws.Cells.Item(20, 3).Value = -2.858730956127765
$ws.Cells.Item(20, 4).Value = 6.065660225717632
$ws.Cells.Item(20, 5).Value = -1.580752749192068
$ws.Cells.Item(20, 6).Value = 3.318631649017334
$ws.Cells.Item(20, 7).Value = -4.592221736907959
$ws.Cells.Item(20, 8).Value = 2.812298059463501

$ws.Cells.Item(21, 3).Value = 3.288362201891221
$ws.Cells.Item(21, 4).Value = -18.18445829341276
$ws.Cells.Item(21, 5).Value = 32.30573925219061
$ws.Cells.Item(21, 6).Value = 6.787132740020752
$ws.Cells.Item(21, 7).Value = 2.726144790649414
$ws.Cells.Item(21, 8).Value = -2.260491132736206

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "walkingToRunning"
$ws.Cells.Item(22, 3).Value = -47.11073684692383
$ws.Cells.Item(22, 4).Value = 3.437827825546265
$ws.Cells.Item(22, 5).Value = -20.08368492126465
$ws.Cells.Item(22, 6).Value = -3.182144641876221
$ws.Cells.Item(22, 7).Value = 2.194577932357788
$ws.Cells.Item(22, 8).Value = -0.2743056118488312

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "walkingToRunning"
$ws.Cells.Item(23, 3).Value = -15.004909515381
$ws.Cells.Item(23, 4).Value = -1.142799942116968
$ws.Cells.Item(23, 5).Value = -13.80460448014111
$ws.Cells.Item(23, 6).Value = -5.860219478607178
$ws.Cells.Item(23, 7).Value = 4.652675628662109
$ws.Cells.Item(23, 8).Value = 1.372593283653259

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "walkingToRunning"
$ws.Cells.Item(24, 3).Value = -7.653207327190201
$ws.Cells.Item(24, 4).Value = -6.902720802708673
$ws.Cells.Item(24, 5).Value = -9.241878660101641
$ws.Cells.Item(24, 6).Value = -3.136338233947754
$ws.Cells.Item(24, 7).Value = 6.668688774108887
$ws.Cells.Item(24, 8).Value = -2.213619470596313

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "walkingToRunning"
$ws.Cells.Item(25, 3).Value = -5.410418410050235
$ws.Cells.Item(25, 4).Value = 5.482645586917334
$ws.Cells.Item(25, 5).Value = 0.0306391966970736
$ws.Cells.Item(25, 6).Value = -0.4337623715400696
$ws.Cells.Item(25, 7).Value = 7.815978527069092
$ws.Cells.Item(25, 8).Value = -6.397445678710938

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "walkingToRunning"
$ws.Cells.Item(26, 3).Value = -5.049034394715962
$ws.Cells.Item(26, 4).Value = 8.959014641611201
$ws.Cells.Item(26, 5).Value = -0.3938065077129187
$ws.Cells.Item(26, 6).Value = 7.690476894378662
$ws.Cells.Item(26, 7).Value = -9.785386085510254
$ws.Cells.Item(26, 8).Value = -1.533448219299316

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "walkingToRunning"
$ws.Cells.Item(27, 3).Value = 9.533257911079737
$ws.Cells.Item(27, 4).Value = -10.20816908384616
$ws.Cells.Item(27, 5).Value = 3.8682978278712
$ws.Cells.Item(27, 6).Value = 1.124719500541687
$ws.Cells.Item(27, 7).Value = -6.268948078155518
$ws.Cells.Item(27, 8).Value = -0.5054680109024048

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "walkingToRunning"
$ws.Cells.Item(28, 3).Value = 18.57404869481136
$ws.Cells.Item(28, 4).Value = -1.860518957439254
$ws.Cells.Item(28, 5).Value = -7.277912340666091
$ws.Cells.Item(28, 6).Value = 0.1537309736013412
$ws.Cells.Item(28, 7).Value = 1.004677534103394
$ws.Cells.Item(28, 8).Value = 0.5731122493743896

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "walkingToRunning"
$ws.Cells.Item(29, 3).Value = -3.77490947121079
$ws.Cells.Item(29, 4).Value = -6.046053083319467
$ws.Cells.Item(29, 5).Value = -15.12075865896126
$ws.Cells.Item(29, 6).Value = -7.232280254364014
$ws.Cells.Item(29, 7).Value = -3.663843870162964
$ws.Cells.Item(29, 8).Value = 11.63588333129883

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "walkingToRunning"
$ws.Cells.Item(30, 3).Value = 3.047089124980793
$ws.Cells.Item(30, 4).Value = -31.00068785014912
$ws.Cells.Item(30, 5).Value = 2.002202786897417
$ws.Cells.Item(30, 6).Value = 1.14389431476593
$ws.Cells.Item(30, 7).Value = 9.762216567993164
$ws.Cells.Item(30, 8).Value = -6.083192825317383

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "walkingToRunning"
$ws.Cells.Item(31, 3).Value = 7.431886748263759
$ws.Cells.Item(31, 4).Value = -15.14435035304032
$ws.Cells.Item(31, 5).Value = 4.078887035972246
$ws.Cells.Item(31, 6).Value = 2.391319036483765
$ws.Cells.Item(31, 7).Value = -6.093712329864502
$ws.Cells.Item(31, 8).Value = 2.177400588989258

